# Updated symbol list: refresh Price, Volume(1h) and Hora columns for each coin row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "304.27"
Set-TextValue $ws.Range("E2") "0.61%"
Set-TextValue $ws.Range("G2") "12"
# Row 3
Set-TextValue $ws.Range("D3") "35.58"
Set-TextValue $ws.Range("E3") "0.90%"
Set-TextValue $ws.Range("G3") "12"
# Row 4
Set-TextValue $ws.Range("D4") "5.060"
Set-TextValue $ws.Range("E4") "0.45%"
Set-TextValue $ws.Range("G4") "12"
# Row 5
Set-TextValue $ws.Range("D5") "0.08061"
Set-TextValue $ws.Range("E5") "0.96%"
Set-TextValue $ws.Range("G5") "12"
# Row 6
Set-TextValue $ws.Range("D6") "1.924"
Set-TextValue $ws.Range("E6") "-0.73%"
Set-TextValue $ws.Range("G6") "12"
# Row 7
Set-TextValue $ws.Range("D7") "4.171"
Set-TextValue $ws.Range("E7") "3.10%"
Set-TextValue $ws.Range("G7") "12"
# Row 8
Set-TextValue $ws.Range("D8") "7.745"
Set-TextValue $ws.Range("E8") "-0.73%"
Set-TextValue $ws.Range("G8") "12"
# Row 9
Set-TextValue $ws.Range("D9") "0.9268"
Set-TextValue $ws.Range("E9") "0.56%"
Set-TextValue $ws.Range("G9") "12"
# Row 10
Set-TextValue $ws.Range("D10") "0.1370"
Set-TextValue $ws.Range("E10") "11.17%"
Set-TextValue $ws.Range("G10") "12"
# Row 11
Set-TextValue $ws.Range("D11") "0.1896"
Set-TextValue $ws.Range("E11") "2.44%"
Set-TextValue $ws.Range("G11") "12"
# Row 12
Set-TextValue $ws.Range("D12") "0.09237"
Set-TextValue $ws.Range("E12") "-3.80%"
Set-TextValue $ws.Range("G12") "12"
# Row 13
Set-TextValue $ws.Range("D13") "0.03572"
Set-TextValue $ws.Range("E13") "1.37%"
Set-TextValue $ws.Range("G13") "12"
# Row 14
Set-TextValue $ws.Range("D14") "0.09811"
Set-TextValue $ws.Range("E14") "-0.43%"
Set-TextValue $ws.Range("G14") "12"
# Row 15
Set-TextValue $ws.Range("D15") "0.001419"
Set-TextValue $ws.Range("E15") "2.14%"
Set-TextValue $ws.Range("G15") "12"
# Row 16
Set-TextValue $ws.Range("D16") "0.005758"
Set-TextValue $ws.Range("E16") "-1.22%"
Set-TextValue $ws.Range("G16") "12"
# Row 17
Set-TextValue $ws.Range("D17") "3.556"
Set-TextValue $ws.Range("E17") "1.43%"
Set-TextValue $ws.Range("G17") "12"
# Row 18
Set-TextValue $ws.Range("D18") "2.896"
Set-TextValue $ws.Range("E18") "-0.84%"
Set-TextValue $ws.Range("G18") "12"
# Row 19
Set-TextValue $ws.Range("E19") "1.98%"
Set-TextValue $ws.Range("G19") "12"
# Row 20
Set-TextValue $ws.Range("D20") "0.1304"
Set-TextValue $ws.Range("E20") "1.10%"
Set-TextValue $ws.Range("G20") "12"
# Row 21
Set-TextValue $ws.Range("D21") "4.892"
Set-TextValue $ws.Range("E21") "-2.63%"
Set-TextValue $ws.Range("G21") "12"
# Row 22
Set-TextValue $ws.Range("D22") "0.2513"
Set-TextValue $ws.Range("E22") "4.97%"
Set-TextValue $ws.Range("G22") "12"
# Row 23
Set-TextValue $ws.Range("D23") "0.04435"
Set-TextValue $ws.Range("E23") "-1.09%"
Set-TextValue $ws.Range("G23") "12"
# Row 24
Set-TextValue $ws.Range("D24") "0.001223"
Set-TextValue $ws.Range("E24") "0.90%"
Set-TextValue $ws.Range("G24") "12"
# Row 25
Set-TextValue $ws.Range("D25") "0.004773"
Set-TextValue $ws.Range("E25") "-0.22%"
Set-TextValue $ws.Range("G25") "12"
# Row 26
Set-TextValue $ws.Range("E26") "32.30%"
Set-TextValue $ws.Range("G26") "12"
# Row 27
Set-TextValue $ws.Range("D27") "0.0003132"
Set-TextValue $ws.Range("E27") "4.61%"
Set-TextValue $ws.Range("G27") "12"
# Row 28
Set-TextValue $ws.Range("G28") "12"
# Row 29
Set-TextValue $ws.Range("G29") "12"
# Row 30
Set-TextValue $ws.Range("G30") "12"
# Row 31
Set-TextValue $ws.Range("G31") "12"
# Row 32
Set-TextValue $ws.Range("G32") "12"
# Row 33
Set-TextValue $ws.Range("G33") "12"
# Row 34
Set-TextValue $ws.Range("G34") "12"
# Row 35
Set-TextValue $ws.Range("G35") "12"
# Row 36
Set-TextValue $ws.Range("G36") "12"
# Row 37
Set-TextValue $ws.Range("G37") "12"
# Row 38
Set-TextValue $ws.Range("G38") "12"
# Row 39
Set-TextValue $ws.Range("D39") "0.01959"
Set-TextValue $ws.Range("E39") "3.07%"
Set-TextValue $ws.Range("G39") "12"
# Row 40
Set-TextValue $ws.Range("D40") "0.04922"
Set-TextValue $ws.Range("E40") "3.93%"
Set-TextValue $ws.Range("G40") "12"
# Row 41
Set-TextValue $ws.Range("D41") "0.007639"
Set-TextValue $ws.Range("E41") "2.60%"
Set-TextValue $ws.Range("G41") "12"
# Row 42
Set-TextValue $ws.Range("E42") "-4.31%"
Set-TextValue $ws.Range("G42") "12"
# Row 43
Set-TextValue $ws.Range("E43") "3.56%"
Set-TextValue $ws.Range("G43") "12"
# Row 44
Set-TextValue $ws.Range("E44") "-0.24%"
Set-TextValue $ws.Range("G44") "12"
# Row 45
Set-TextValue $ws.Range("D45") "0.01077"
Set-TextValue $ws.Range("E45") "-0.11%"
Set-TextValue $ws.Range("G45") "12"
# Row 46
Set-TextValue $ws.Range("D46") "0.00006372"
Set-TextValue $ws.Range("E46") "2.31%"
Set-TextValue $ws.Range("G46") "12"
# Row 47
Set-TextValue $ws.Range("E47") "0.23%"
Set-TextValue $ws.Range("G47") "12"
# Row 48
Set-TextValue $ws.Range("D48") "64.96"
Set-TextValue $ws.Range("E48") "0.75%"
Set-TextValue $ws.Range("G48") "12"
# Row 49
Set-TextValue $ws.Range("D49") "0.001192"
Set-TextValue $ws.Range("E49") "-19.99%"
Set-TextValue $ws.Range("G49") "12"
# Row 50
Set-TextValue $ws.Range("E50") "0.23%"
Set-TextValue $ws.Range("G50") "12"
# Row 51
Set-TextValue $ws.Range("E51") "0.23%"
Set-TextValue $ws.Range("G51") "12"
